$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G6").Value = "D4*"
Write-Host "G6:" $ws.Range("G6").Value
